# Shift the daily schedule (Sheet1) forward by 20 minutes for the
# wake-up slot and for the evening block starting at row 13.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 05:30 -> 05:50 (wake up start time)
$ws.Range("B1").Value = 0.24305555555555555
# 05:50 -> 06:00 (wake up end time / next slot start)
$ws.Range("C1").Value = 0.25

# 20:10 -> 20:30
$ws.Range("C13").Value = 0.85416666666666663
# 22:00 -> 22:20
$ws.Range("C14").Value = 0.93055555555555547
# 22:30 -> 22:50
$ws.Range("C15").Value = 0.95138888888888884

# Move the active selection, matching the workbook's recorded view state.
$ws.Range("C15").Select()
